$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2105.348
$ws.Range("I19").Value = 2091.1
$ws.Range("J19").Value = 2116.3076
$ws.Range("K19").Value = 2091.1
$ws.Range("L19").Value = 2116.3076
$ws.Range("M19").Value = -1916.1
$ws.Range("N19").Value = -2466.3076

$ws.Range("H43").Value = 8752.6
$ws.Range("J43").Value = 8691
$ws.Range("L43").Value = 8691
$ws.Range("N43").Value = -8829

$ws.Range("H74").Value = 5500.6
$ws.Range("I74").Value = 4876
$ws.Range("K74").Value = 4876
$ws.Range("M74").Value = -3940

$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 5000
$ws.Range("N76").Value = -5630
$ws.Range("M76").ClearContents()

$ws.Range("H77").Value = 5500.6
$ws.Range("I77").Value = 4876
$ws.Range("K77").Value = 24380
$ws.Range("M77").Value = -19700

$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 5000
$ws.Range("N79").Value = -7184
$ws.Range("M79").ClearContents()

$ws.Range("H80").Value = 3448.4546
$ws.Range("I80").Value = 1604.1818
$ws.Range("K80").Value = 4812.5454
$ws.Range("M80").Value = -3814.5454

$ws.Range("H83").Value = 3448.4546
$ws.Range("I83").Value = 1604.1818
$ws.Range("K83").Value = 14437.6362
$ws.Range("M83").Value = -9445.636200000001

$ws.Range("H99").Value = 1456.9
$ws.Range("I99").Value = 267
$ws.Range("K99").Value = 801
$ws.Range("M99").Value = 697

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H116").Value = 150599.58
$ws.Range("J116").Value = 9033
$ws.Range("L116").Value = 9033
$ws.Range("N116").Value = -15917

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 96295.664
$ws.Range("J52").Value = 96295.664
$ws.Range("L52").Value = 96295.664
$ws.Range("N52").Value = -96931.664

$ws.Range("H63").Value = 1890.375
$ws.Range("I63").Value = 1768
$ws.Range("K63").Value = 1768
$ws.Range("M63").Value = -1082

$ws.Range("H66").Value = 1890.375
$ws.Range("I66").Value = 1768
$ws.Range("K66").Value = 8840
$ws.Range("M66").Value = -5408

$ws.Range("H97").Value = 1644.6
$ws.Range("I97").Value = 761.1667
$ws.Range("K97").Value = 761.1667
$ws.Range("M97").Value = -265.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2220.6155
$ws.Range("I105").Value = 2137
$ws.Range("J105").Value = 2499.3333
$ws.Range("K105").Value = 2137
$ws.Range("L105").Value = 2499.3333
$ws.Range("M105").Value = -390
$ws.Range("N105").Value = -5993.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3453.5
$ws.Range("I31").Value = 1933.5
$ws.Range("K31").Value = 1933.5
$ws.Range("M31").Value = -1638.5

$ws.Range("H34").Value = 3453.5
$ws.Range("I34").Value = 1933.5
$ws.Range("K34").Value = 1933.5
$ws.Range("M34").Value = -1731.5

$ws.Range("H51").Value = 33924.188
$ws.Range("I51").Value = 13128
$ws.Range("K51").Value = 13128
$ws.Range("M51").Value = -12392

$ws.Range("H61").Value = 33924.188
$ws.Range("I61").Value = 13128
$ws.Range("K61").Value = 13128
$ws.Range("M61").Value = -12780

$ws.Range("H94").Value = 610.94116
$ws.Range("I94").Value = 966
$ws.Range("J94").Value = 563.6
$ws.Range("K94").Value = 966
$ws.Range("L94").Value = 563.6
$ws.Range("M94").Value = -515
$ws.Range("N94").Value = -1465.6

$ws.Range("H99").Value = 2219
$ws.Range("I99").Value = 2168.8
$ws.Range("K99").Value = 2168.8
$ws.Range("M99").Value = -670.8000000000002

$ws.Range("H105").Value = 2886.625
$ws.Range("I105").Value = 2886.625
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2886.625
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1139.625
$ws.Range("N105").ClearContents()

$ws.Range("H126").Value = 2219
$ws.Range("I126").Value = 2168.8
$ws.Range("K126").Value = 6506.400000000001
$ws.Range("M126").Value = -4036.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 4010.8667
$ws.Range("I11").Value = 5266.5557
$ws.Range("K11").Value = 15799.6671
$ws.Range("M11").Value = -15659.6671

$ws.Range("H39").Value = 6324.476
$ws.Range("I39").Value = 3650.5
$ws.Range("J39").Value = 7394.067
$ws.Range("K39").Value = 10951.5
$ws.Range("L39").Value = 22182.201
$ws.Range("M39").Value = -10657.5
$ws.Range("N39").Value = -22770.201

$ws.Range("H55").Value = 4536.625
$ws.Range("I55").Value = 901.5
$ws.Range("J55").Value = 5748.3335
$ws.Range("K55").Value = 2704.5
$ws.Range("L55").Value = 17245.0005
$ws.Range("M55").Value = -2527.5
$ws.Range("N55").Value = -17599.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3677442.5
$ws.Range("I11").Value = 3929553.5
$ws.Range("J11").Value = 400000
$ws.Range("K11").Value = 3929553.5
$ws.Range("L11").Value = 400000
$ws.Range("M11").Value = -3929414.5
$ws.Range("N11").Value = -400278

$ws.Range("H15").Value = 17118
$ws.Range("J15").Value = 17118
$ws.Range("L15").Value = 17118
$ws.Range("N15").Value = -17694

$ws.Range("H81").Value = 17118
$ws.Range("J81").Value = 17118
$ws.Range("L81").Value = 17118
$ws.Range("N81").Value = -19114

$ws.Range("H84").Value = 17118
$ws.Range("J84").Value = 17118
$ws.Range("L84").Value = 51354
$ws.Range("N84").Value = -61338

$ws.Range("H102").Value = 3041.6428
$ws.Range("I102").Value = 2328.3333
$ws.Range("J102").Value = 3236.182
$ws.Range("K102").Value = 2328.3333
$ws.Range("L102").Value = 3236.182
$ws.Range("M102").Value = -706.3332999999998
$ws.Range("N102").Value = -6480.182

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 2351.1667
$ws.Range("I132").Value = 2009.4546
$ws.Range("J132").Value = 2888.1428
$ws.Range("K132").Value = 6028.3638
$ws.Range("L132").Value = 8664.428400000001
$ws.Range("M132").Value = -3498.3638
$ws.Range("N132").Value = -13724.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 165000
$ws.Range("J24").Value = 165000
$ws.Range("L24").Value = 165000
$ws.Range("N24").Value = -165686

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 37943.777
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H100").Value = 2718.0625
$ws.Range("I100").Value = 2187
$ws.Range("J100").Value = 3249.125
$ws.Range("K100").Value = 4374
$ws.Range("L100").Value = 6498.25
$ws.Range("M100").Value = -3833
$ws.Range("N100").Value = -7580.25
